# Weekly update: insert a new price record as row 130 for
# "Feria Lagunitas de Puerto Montt - Mango", pushing the existing
# rows 130-150 down to 131-151.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 130:150 down one row, creating a blank row 130.
$ws.Rows("130:130").Insert()

# Populate the new row 130 with this week's record.
$ws.Range("A130").Value = 4
$ws.Range("B130").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C130").Value = "Los Lagos"
$ws.Range("D130").Value = 44578
$ws.Range("E130").Value = 10
$ws.Range("F130").Value = "Fruta"
$ws.Range("G130").Value = 100108
$ws.Range("H130").Value = "Tropicales y subtropicales"
$ws.Range("I130").Value = 100108002
$ws.Range("J130").Value = "Mango"
$ws.Range("K130").Value = "Sin especificar"
$ws.Range("L130").Value = "Primera"
$ws.Range("M130").Value = 80
$ws.Range("N130").Value = 8000
$ws.Range("O130").Value = 8500
$ws.Range("P130").Value = 8250
$ws.Range("Q130").Value = "`$/bandeja 4 kilos"
$ws.Range("R130").Value = "Perú"
$ws.Range("S130").Value = 2062
$ws.Range("T130").Value = 4
